$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": add a new day column DD (29-sep) after DC (28-sep)
# ---------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing header cell (DC1) into the
# new header cell (DD1) so it keeps the bold / centered / bordered style.
$wsSpot.Range("DC1").Copy($wsSpot.Range("DD1"))
$wsSpot.Range("DD1").Value = "29-sep"

$wsSpot.Range("DD2").Value = 51.6
$wsSpot.Range("DD3").Value = 47.7
$wsSpot.Range("DD4").Value = 45
$wsSpot.Range("DD5").Value = 32
$wsSpot.Range("DD6").Value = 30
$wsSpot.Range("DD7").Value = 32.5
$wsSpot.Range("DD8").Value = 50
$wsSpot.Range("DD9").Value = 74.40000000000001
$wsSpot.Range("DD10").Value = 81.8
$wsSpot.Range("DD11").Value = 73.08
$wsSpot.Range("DD12").Value = 59.33
$wsSpot.Range("DD13").Value = 40.14
$wsSpot.Range("DD14").Value = 35
$wsSpot.Range("DD15").Value = 28.3
$wsSpot.Range("DD16").Value = 21.99
$wsSpot.Range("DD17").Value = 27.94
$wsSpot.Range("DD18").Value = 35
$wsSpot.Range("DD19").Value = 52
$wsSpot.Range("DD20").Value = 85
$wsSpot.Range("DD21").Value = 101.13
$wsSpot.Range("DD22").Value = 133.26
$wsSpot.Range("DD23").Value = 89.51000000000001
$wsSpot.Range("DD24").Value = 84.88
$wsSpot.Range("DD25").Value = 88.59999999999999

# ---------------------------------------------------------------
# Sheet "Gaz": append two new daily rows (105, 106)
# ---------------------------------------------------------------
# The date-like strings in column A must stay plain text (as in the
# rest of the column) instead of being auto-converted to a date
# serial number, so the cell is briefly marked as Text ("@"), the
# literal value is written, and the temporary formatting is cleared
# again so the cell keeps the same (unstyled) look as its neighbours.
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A105").NumberFormat = "@"
$wsGaz.Range("A105").Value = "2025-09-27"
$wsGaz.Range("A105").ClearFormats()
$wsGaz.Range("B105").Value = 31.775

$wsGaz.Range("A106").NumberFormat = "@"
$wsGaz.Range("A106").Value = "2025-09-28"
$wsGaz.Range("A106").ClearFormats()
$wsGaz.Range("B106").Value = 31.775

# ---------------------------------------------------------------
# Sheet "CO2": append two new daily rows (105, 106)
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A105").NumberFormat = "@"
$wsCo2.Range("A105").Value = "2025-09-27"
$wsCo2.Range("A105").ClearFormats()
$wsCo2.Range("B105").Value = 75.26000000000001

$wsCo2.Range("A106").NumberFormat = "@"
$wsCo2.Range("A106").Value = "2025-09-28"
$wsCo2.Range("A106").ClearFormats()
$wsCo2.Range("B106").Value = 75.26000000000001
